$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 237, shifting existing rows 237-328 down to 238-329.
$ws.Rows("237:237").Insert()

# Populate the newly inserted row 237 with the new record.
$ws.Range("A237").Value = 5
$ws.Range("B237").Value = "Macroferia Regional de Talca"
$ws.Range("C237").Value = "Maule"
$ws.Range("D237").Value = 45119
$ws.Range("E237").Value = 7
$ws.Range("F237").Value = 100112017
$ws.Range("G237").Value = "Apio"
$ws.Range("H237").Value = "Americana (o)"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 700
$ws.Range("K237").Value = 6000
$ws.Range("L237").Value = 6000
$ws.Range("M237").Value = 6000
$ws.Range("N237").Value = "$/docena de matas"
$ws.Range("O237").Value = "Provincia del Elquí"
$ws.Range("P237").Value = 1000
$ws.Range("Q237").Value = 6
$ws.Range("R237").Value = "Hortaliza"
